$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MASSO")
$ws.Select()

# Clear the shared formulas in B4:B8 and C4:C8 so each cell becomes a literal value
$ws.Range("B4:B8").ClearContents()
$ws.Range("C4:C8").ClearContents()

# Set the new availability text, in the order that reproduces the target
# shared-string table ordering.
$ws.Range("B3").Value = "20h à 22h - Claudie Germain"
$ws.Range("C3").Value = "8PM to 10PM - Claudie Germain"
$ws.Range("C9").Value = "6PM to 8PM - Claudie Germain"
$ws.Range("B6").Value = "7h30 à 9h30 - Marie Pier Pépin"
$ws.Range("B8").Value = "18h à 20h - Claudie Germain"
$ws.Range("C6").Value = "7:30AM tp  9:30AM - Marie Pier Pépin"
$ws.Range("B2").Value = "Pas de disponibilité"
$ws.Range("C2").Value = "No availability"

$ws.Range("B4").Value = "Pas de disponibilité"
$ws.Range("C4").Value = "No availability"
$ws.Range("B5").Value = "Pas de disponibilité"
$ws.Range("C5").Value = "No availability"
$ws.Range("B7").Value = "Pas de disponibilité"
$ws.Range("C7").Value = "No availability"
$ws.Range("C8").Value = "No availability"

$ws.Range("B15").Select()
